$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 32 (RLC): turn it into the "Camera ID" row ---
# B32 becomes the new shared string "Camera ID" (added first so it lands at
# the lower shared-string index, matching the upstream edit order) and gets
# the green "needs definition" highlight (matches the fill used elsewhere,
# e.g. B3 originally / B30 / B31).
$ws.Range("B32").Value = "Camera ID"
$ws.Range("B32").Interior.Color = 5296274

# C32's "X" marker is removed (no longer applicable).
$ws.Range("C32").ClearContents()

# --- Row 3 (INTERSECTION_ID): add "Camera ID"-style annotation columns ---
# C3 gets the "X" marker (reuses existing shared string "X")
$ws.Range("C3").Value = "X"

# D3 gets the new Notes text (new shared string)
$ws.Range("D3").Value = "Has null values.  Not required.  RLC column as unique id for red light camera"

# B3 loses its "undefined/needs formatting" highlight (was s=3 with a green
# fill); copy the plain bordered, non-highlighted format from a cell that
# already carries that look (B25) onto B3, without touching its text.
$ws.Range("B25").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# Row 3 grows taller to fit the wrapped note text in D3.
$ws.Rows.Item(3).RowHeight = 30

# --- Selection moves to A3 ---
$ws.Range("A3").Select()
